$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing cell contents in the used range (A1:K7). Styles on
# individual cells are preserved by ClearContents, only the values/shared
# string usages are dropped. Re-writing everything from scratch (in the
# exact order below) lets the shared-strings table come out in the same
# order as the target file: existing header labels first (their original
# order), then the three brand-new header labels, then the row labels
# (aon/flex/sub) which are "re-discovered" after the new ones.
$ws.Range("A1:K7").ClearContents()

# --- Row 1: header -------------------------------------------------------
$ws.Range("A1").Value = "geral_modalidade"
$ws.Range("B1").Value = "mencoes_ccxp"
$ws.Range("C1").Value = "total"
$ws.Range("D1").Value = "total_sucesso"
$ws.Range("E1").Value = "particip"
$ws.Range("F1").Value = "taxa_sucesso"
$ws.Range("G1").Value = "arrecadado_sucesso"
$ws.Range("H1").Value = "media_sucesso"
$ws.Range("I1").Value = "std_sucesso"
$ws.Range("J1").Value = "min_sucesso"
$ws.Range("K1").Value = "max_sucesso"

# New header columns (apply the same header style as the existing headers
# by copying K1's formatting - ClearContents keeps K1's original style, so
# this reproduces the bold / centered / bordered header look on L1:N1)
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 -----------------------------------------------------------------
$ws.Range("A2").Value = "aon"
$ws.Range("B2").Value = $false
$ws.Range("C2").Value = 1179
$ws.Range("D2").Value = 704
$ws.Range("E2").Value = 88.31460674157303
$ws.Range("F2").Value = 59.71162001696353
$ws.Range("G2").Value = 21202461.38872604
$ws.Range("H2").Value = 30117.1326544404
$ws.Range("I2").Value = 47960.70656678714
$ws.Range("J2").Value = 41.81688448509265
$ws.Range("K2").Value = 679297.6600721752
$ws.Range("L2").Value = 92.57989061438856
$ws.Range("M2").Value = 229018
$ws.Range("N2").Value = 325.3096590909091

# --- Row 3 -----------------------------------------------------------------
$ws.Range("A3").Value = "aon"
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = 156
$ws.Range("D3").Value = 126
$ws.Range("E3").Value = 11.68539325842697
$ws.Range("F3").Value = 80.76923076923077
$ws.Range("G3").Value = 2860818.438596986
$ws.Range("H3").Value = 22704.90824283323
$ws.Range("I3").Value = 20582.18081328499
$ws.Range("J3").Value = 1720.659275370021
$ws.Range("K3").Value = 154365.9837040891
$ws.Range("L3").Value = 82.83823479360029
$ws.Range("M3").Value = 34535
$ws.Range("N3").Value = 274.0873015873016

# --- Row 4 -----------------------------------------------------------------
$ws.Range("A4").Value = "flex"
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 1293
$ws.Range("D4").Value = 1208
$ws.Range("E4").Value = 88.07901907356948
$ws.Range("F4").Value = 93.4261407579273
$ws.Range("G4").Value = 15952537.73508588
$ws.Range("H4").Value = 13205.74315818367
$ws.Range("I4").Value = 35396.12153793828
$ws.Range("J4").Value = 10.77163914429046
$ws.Range("K4").Value = 708972.7845446636
$ws.Range("L4").Value = 89.86130176813208
$ws.Range("M4").Value = 177524
$ws.Range("N4").Value = 146.9569536423841

# --- Row 5 -----------------------------------------------------------------
$ws.Range("A5").Value = "flex"
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = 175
$ws.Range("D5").Value = 175
$ws.Range("E5").Value = 11.92098092643052
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 2409594.202473221
$ws.Range("H5").Value = 13769.10972841841
$ws.Range("I5").Value = 21333.70082549344
$ws.Range("J5").Value = 313.2716886535385
$ws.Range("K5").Value = 121747.7984910872
$ws.Range("L5").Value = 92.24386350483199
$ws.Range("M5").Value = 26122
$ws.Range("N5").Value = 149.2685714285714

# --- Row 6 -----------------------------------------------------------------
$ws.Range("A6").Value = "sub"
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 678
$ws.Range("D6").Value = 147
$ws.Range("E6").Value = 99.12280701754386
$ws.Range("F6").Value = 21.68141592920354
$ws.Range("G6").Value = 41148.96897529707
$ws.Range("H6").Value = 279.9249590156264
$ws.Range("I6").Value = 649.3661034486928
$ws.Range("J6").Value = 1.087396962410123
$ws.Range("K6").Value = 5087.076865717208
$ws.Range("L6").Value = 19.47419260544111
$ws.Range("M6").Value = 2113
$ws.Range("N6").Value = 14.37414965986395

# --- Row 7 -----------------------------------------------------------------
$ws.Range("A7").Value = "sub"
$ws.Range("B7").Value = $true
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 0.8771929824561403
$ws.Range("F7").Value = 83.33333333333334
$ws.Range("G7").Value = 2037.988779487728
$ws.Range("H7").Value = 407.5977558975457
$ws.Range("I7").Value = 752.9936319265861
$ws.Range("J7").Value = 40.66419228170764
$ws.Range("K7").Value = 1753.365733305352
$ws.Range("L7").Value = 21.45251346829188
$ws.Range("M7").Value = 95
$ws.Range("N7").Value = 19
